$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1869.7667
$ws.Range("I40").Value = 1828.2609
$ws.Range("J40").Value = 2006.1428
$ws.Range("K40").Value = 1828.2609
$ws.Range("L40").Value = 2006.1428
$ws.Range("M40").Value = -1653.2609
$ws.Range("N40").Value = -2356.1428
# Row 43
$ws.Range("H43").Value = 1100
$ws.Range("J43").Value = 1233.3334
$ws.Range("L43").Value = 1233.3334
$ws.Range("N43").Value = -1371.3334
# Row 115
$ws.Range("H115").Value = 728
$ws.Range("I115").Value = 728
$ws.Range("K115").Value = 2184
$ws.Range("M115").Value = -617
# Row 138
$ws.Range("H138").Value = 5329.9795
$ws.Range("I138").Value = 909.03845
$ws.Range("J138").Value = 10327.565
$ws.Range("K138").Value = 2727.11535
$ws.Range("L138").Value = 30982.695
$ws.Range("M138").Value = 2412.88465
$ws.Range("N138").Value = -41262.695

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3501.3333
$ws.Range("I2").Value = 3358.7144
$ws.Range("J2").Value = 4000.5
$ws.Range("K2").Value = 3358.7144
$ws.Range("L2").Value = 4000.5
$ws.Range("M2").Value = -3245.7144
$ws.Range("N2").Value = -4226.5
# Row 45
$ws.Range("H45").Value = 5174.76
$ws.Range("I45").Value = 7324.9375
$ws.Range("J45").Value = 1352.2222
$ws.Range("K45").Value = 7324.9375
$ws.Range("L45").Value = 1352.2222
$ws.Range("M45").Value = -6947.9375
$ws.Range("N45").Value = -2106.2222
# Row 61
$ws.Range("H61").Value = 5886.44
$ws.Range("I61").Value = 6228.7393
$ws.Range("K61").Value = 6228.7393
$ws.Range("M61").Value = -6016.7393
# Row 63
$ws.Range("H63").Value = 142860270
$ws.Range("I63").Value = 166669490
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 166669490
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -166668804
$ws.Range("N63").Value = -6372
# Row 66
$ws.Range("H66").Value = 142860270
$ws.Range("I66").Value = 166669490
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 833347450
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -833344018
$ws.Range("N66").Value = -31864
# Row 74
$ws.Range("H74").Value = 1825.5098
$ws.Range("I74").Value = 1734.0638
$ws.Range("J74").Value = 2900
$ws.Range("K74").Value = 1734.0638
$ws.Range("L74").Value = 2900
$ws.Range("M74").Value = -860.0637999999999
$ws.Range("N74").Value = -4648
# Row 77
$ws.Range("H77").Value = 1825.5098
$ws.Range("I77").Value = 1734.0638
$ws.Range("J77").Value = 2900
$ws.Range("K77").Value = 8670.319
$ws.Range("L77").Value = 14500
$ws.Range("M77").Value = -4302.319
$ws.Range("N77").Value = -23236
# Row 116
$ws.Range("H116").Value = 3501.3333
$ws.Range("I116").Value = 3358.7144
$ws.Range("J116").Value = 4000.5
$ws.Range("K116").Value = 3358.7144
$ws.Range("L116").Value = 4000.5
$ws.Range("M116").Value = -1064.7144
$ws.Range("N116").Value = -8588.5
# Row 122
$ws.Range("H122").Value = 1352788.9
$ws.Range("I122").Value = 1711729.9
$ws.Range("K122").Value = 5135189.699999999
$ws.Range("M122").Value = -5132739.699999999
# Row 132
$ws.Range("H132").Value = 3472.225
$ws.Range("I132").Value = 1782.0952
$ws.Range("J132").Value = 5340.263
$ws.Range("K132").Value = 5346.2856
$ws.Range("L132").Value = 16020.789
$ws.Range("M132").Value = -2816.2856
$ws.Range("N132").Value = -21080.789
# Row 136
$ws.Range("H136").Value = 5886.44
$ws.Range("I136").Value = 6228.7393
$ws.Range("K136").Value = 18686.2179
$ws.Range("M136").Value = -16136.2179

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3501.3333
$ws.Range("I3").Value = 3358.7144
$ws.Range("J3").Value = 4000.5
$ws.Range("K3").Value = 3358.7144
$ws.Range("L3").Value = 4000.5
$ws.Range("M3").Value = -3244.7144
$ws.Range("N3").Value = -4228.5
# Row 55
$ws.Range("H55").Value = 79800
$ws.Range("J55").Value = 79800
$ws.Range("L55").Value = 79800
$ws.Range("N55").Value = -80346
# Row 107
$ws.Range("H107").Value = 829.1429000000001
$ws.Range("I107").Value = 823.3333
$ws.Range("J107").Value = 846.5714
$ws.Range("K107").Value = 823.3333
$ws.Range("L107").Value = 846.5714
$ws.Range("M107").Value = 1096.6667
$ws.Range("N107").Value = -4686.5714

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 6411513.5
$ws.Range("I16").Value = 10990037
$ws.Range("J16").Value = 1580
$ws.Range("K16").Value = 10990037
$ws.Range("L16").Value = 1580
$ws.Range("M16").Value = -10989750
$ws.Range("N16").Value = -2154
# Row 31
$ws.Range("H31").Value = 7458.9565
$ws.Range("I31").Value = 1795.8667
$ws.Range("K31").Value = 1795.8667
$ws.Range("M31").Value = -1500.8667
# Row 34
$ws.Range("H34").Value = 7458.9565
$ws.Range("I34").Value = 1795.8667
$ws.Range("K34").Value = 1795.8667
$ws.Range("M34").Value = -1593.8667
# Row 58
$ws.Range("H58").Value = 1554.1708
$ws.Range("I58").Value = 776.9545000000001
$ws.Range("J58").Value = 2454.1052
$ws.Range("K58").Value = 776.9545000000001
$ws.Range("L58").Value = 2454.1052
$ws.Range("M58").Value = -573.9545000000001
$ws.Range("N58").Value = -2860.1052
# Row 86
$ws.Range("H86").Value = 2323.4
$ws.Range("I86").Value = 2313.8333
$ws.Range("J86").Value = 2361.6667
$ws.Range("K86").Value = 2313.8333
$ws.Range("L86").Value = 2361.6667
$ws.Range("M86").Value = -1190.8333
$ws.Range("N86").Value = -4607.6667
# Row 89
$ws.Range("H89").Value = 2323.4
$ws.Range("I89").Value = 2313.8333
$ws.Range("J89").Value = 2361.6667
$ws.Range("K89").Value = 11569.1665
$ws.Range("L89").Value = 11808.3335
$ws.Range("M89").Value = -5953.166499999999
$ws.Range("N89").Value = -23040.3335
# Row 94
$ws.Range("H94").Value = 2730.4517
$ws.Range("J94").Value = 2293.9565
$ws.Range("L94").Value = 2293.9565
$ws.Range("N94").Value = -3195.9565
# Row 99
$ws.Range("H99").Value = 4312916
$ws.Range("I99").Value = 1635.8948
$ws.Range("J99").Value = 12504349
$ws.Range("K99").Value = 1635.8948
$ws.Range("L99").Value = 12504349
$ws.Range("M99").Value = -137.8948
$ws.Range("N99").Value = -12507345
# Row 105
$ws.Range("H105").Value = 33335176
$ws.Range("I105").Value = 55558092
$ws.Range("K105").Value = 55558092
$ws.Range("M105").Value = -55556345
# Row 113
$ws.Range("H113").Value = 6411513.5
$ws.Range("I113").Value = 10990037
$ws.Range("J113").Value = 1580
$ws.Range("K113").Value = 10990037
$ws.Range("L113").Value = 1580
$ws.Range("M113").Value = -10987867
$ws.Range("N113").Value = -5920
# Row 122
$ws.Range("H122").Value = 1244.0454
$ws.Range("I122").Value = 1086.5883
$ws.Range("K122").Value = 3259.7649
$ws.Range("M122").Value = -809.7648999999997
# Row 126
$ws.Range("H126").Value = 4312916
$ws.Range("I126").Value = 1635.8948
$ws.Range("J126").Value = 12504349
$ws.Range("K126").Value = 4907.6844
$ws.Range("L126").Value = 37513047
$ws.Range("M126").Value = -2437.6844
$ws.Range("N126").Value = -37517987
# Row 132
$ws.Range("H132").Value = 2765.2778
$ws.Range("I132").Value = 2305.923
$ws.Range("J132").Value = 3959.6
$ws.Range("K132").Value = 6917.768999999999
$ws.Range("L132").Value = 11878.8
$ws.Range("M132").Value = -4387.768999999999
$ws.Range("N132").Value = -16938.8
# Row 134
$ws.Range("H134").Value = 4715.8125
$ws.Range("I134").Value = 7025
$ws.Range("K134").Value = 21075
$ws.Range("M134").Value = -18540
# Row 136
$ws.Range("H136").Value = 1554.1708
$ws.Range("I136").Value = 776.9545000000001
$ws.Range("J136").Value = 2454.1052
$ws.Range("K136").Value = 2330.8635
$ws.Range("L136").Value = 7362.3156
$ws.Range("M136").Value = 219.1364999999996
$ws.Range("N136").Value = -12462.3156

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 142858740
$ws.Range("I113").Value = 500000500
$ws.Range("J113").Value = 2038.6
$ws.Range("K113").Value = 500000500
$ws.Range("L113").Value = 2038.6
$ws.Range("M113").Value = -499998330
$ws.Range("N113").Value = -6378.6
# Row 126
$ws.Range("H126").Value = 6965.8
$ws.Range("I126").Value = 9516.538
$ws.Range("J126").Value = 2228.7144
$ws.Range("K126").Value = 28549.614
$ws.Range("L126").Value = 6686.1432
$ws.Range("M126").Value = -26079.614
$ws.Range("N126").Value = -11626.1432
# Row 132
$ws.Range("H132").Value = 5315.8335
$ws.Range("I132").Value = 13369.75
$ws.Range("J132").Value = 3014.7144
$ws.Range("K132").Value = 40109.25
$ws.Range("L132").Value = 9044.143199999999
$ws.Range("M132").Value = -37579.25
$ws.Range("N132").Value = -14104.1432

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1573.875
$ws.Range("I126").Value = 1013
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 3039
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -569
$ws.Range("N126").Value = -21440
# Row 132
$ws.Range("H132").Value = 2653.6667
$ws.Range("I132").Value = 1983.5
$ws.Range("J132").Value = 3323.8333
$ws.Range("K132").Value = 5950.5
$ws.Range("L132").Value = 9971.499899999999
$ws.Range("M132").Value = -3420.5
$ws.Range("N132").Value = -15031.4999
